$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B91").Value = 0.01624
$ws.Range("B92").Value = 0.01831
$ws.Range("B93").Value = 0.01857
$ws.Range("B95").Value = 0.02898
$ws.Range("B96").Value = 0.02634
$ws.Range("B97").Value = 0.01848
$ws.Range("B98").Value = 0.04465
$ws.Range("B99").Value = 0.02256
$ws.Range("B100").Value = 0.02488
$ws.Range("B102").Value = 0.01939
$ws.Range("B103").Value = 0.01448
$ws.Range("B104").Value = 0.02089
$ws.Range("B105").Value = 0.03078
$ws.Range("B106").Value = 0.03055
$ws.Range("B107").Value = 0.01378
$ws.Range("B108").Value = 0.01976
$ws.Range("B109").Value = 0.01314
$ws.Range("B110").Value = 0.02428
$ws.Range("B111").Value = 0.02579
$ws.Range("B135").Value = 0.01716
$ws.Range("B136").Value = 0.00516
$ws.Range("B137").Value = 0.01298
$ws.Range("B138").Value = 0.01321
$ws.Range("B139").Value = 0.01485
$ws.Range("B140").Value = 0.01327
$ws.Range("B141").Value = 0.01554
$ws.Range("B142").Value = 0.01384
$ws.Range("B143").Value = 0.00979
$ws.Range("B144").Value = 0.02713
$ws.Range("B146").Value = 0.01524
$ws.Range("B147").Value = 0.00916
$ws.Range("B148").Value = 0.01068
$ws.Range("B149").Value = 0.02356
$ws.Range("B150").Value = 0.03857
$ws.Range("B151").Value = 0.01037
$ws.Range("B152").Value = 0.01151
$ws.Range("B153").Value = 0.01089
$ws.Range("B154").Value = 0.00782
$ws.Range("B155").Value = 0.02239
$ws.Range("B168").Value = 0.01545
$ws.Range("B169").Value = 0.02702
$ws.Range("B170").Value = 0.01115
$ws.Range("B171").Value = 0.0118
$ws.Range("B172").Value = 0.01357
$ws.Range("B173").Value = 0.00845
$ws.Range("B174").Value = 0.01215
$ws.Range("B175").Value = 0.01008
$ws.Range("B176").Value = 0.00973
$ws.Range("B177").Value = 0.02288

# Update sheet view: scroll position and active selection
$ws.Activate()
$ws.Range("C143").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 132
$win.ScrollColumn = 1
